$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns store values as text (e.g. "19.10", "48.50") even
# though some look numeric; force text format first so Excel's auto-detect
# doesn't coerce them into numbers and strip significant trailing zeros.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "44.436.29"
$ws.Range("E2").Value = "  +3.45%  "
$ws.Range("D3").Value = "2.417.70"
$ws.Range("E3").Value = "  +2.31%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "313.75"
$ws.Range("E5").Value = "  +3.65%  "
$ws.Range("D6").Value = "100.74"
$ws.Range("E6").Value = "  +5.02%  "
$ws.Range("E7").Value = "  +2.00%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "0.511"
$ws.Range("E9").Value = "  +4.33%  "
$ws.Range("D10").Value = "35.14"
$ws.Range("E10").Value = "  +2.95%  "
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("D12").Value = "19.10"
$ws.Range("E12").Value = "  +4.34%  "
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("E14").Value = "  +2.94%  "
$ws.Range("D15").Value = "2.797.12"
$ws.Range("E15").Value = "  +2.43%  "
$ws.Range("D16").Value = "2.432.43"
$ws.Range("E16").Value = "  +3.26%  "
$ws.Range("D17").Value = "0.831"
$ws.Range("E17").Value = "  +4.89%  "
$ws.Range("D18").Value = "44.308.19"
$ws.Range("E18").Value = "  +3.22%  "
$ws.Range("D19").Value = "12.38"
$ws.Range("E19").Value = "  +4.26%  "
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").Value = "0.0₃0919"
$ws.Range("E21").Value = "  +3.94%  "
$ws.Range("D22").Value = "68.59"
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("D23").Value = "241.51"
$ws.Range("E23").Value = "  +2.81%  "
$ws.Range("D24").Value = "2.27"
$ws.Range("E24").Value = "  +5.36%  "
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "25.08"
$ws.Range("E27").Value = "  +2.55%  "
$ws.Range("D28").Value = "2.28"
$ws.Range("E28").Value = "  -4.07%  "
$ws.Range("E29").Value = "  +2.41%  "
$ws.Range("D30").Value = "33.15"
$ws.Range("E30").Value = "  +3.46%  "
$ws.Range("D31").Value = "48.50"
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("E32").Value = "  +18.78%  "
$ws.Range("D33").Value = "19.28"
$ws.Range("E33").Value = "  +11.13%  "
$ws.Range("D34").Value = "5.16"
$ws.Range("E34").Value = "  +2.76%  "
$ws.Range("D35").Value = "0.0771"
$ws.Range("E35").Value = "  +7.62%  "
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("D37").Value = "1.88"
$ws.Range("E37").Value = "  +1.78%  "
$ws.Range("D38").Value = "4.48"
$ws.Range("E38").Value = "  +3.02%  "
$ws.Range("E39").Value = "  +2.71%  "
$ws.Range("D40").Value = "122.57"
$ws.Range("E40").Value = "  -5.64%  "
$ws.Range("E41").Value = "  -2.30%  "
$ws.Range("E42").Value = "  +1.40%  "
$ws.Range("D43").Value = "21.17"
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("E44").Value = "  +3.52%  "
$ws.Range("D45").Value = "1.946.03"
$ws.Range("E45").Value = "  +0.81%  "
$ws.Range("E46").Value = "  +2.23%  "
$ws.Range("D47").Value = "2.92"
$ws.Range("E47").Value = "  +7.69%  "
$ws.Range("D48").Value = "9.48"
$ws.Range("E48").Value = "  +3.67%  "
$ws.Range("E49").Value = "  +9.58%  "
$ws.Range("D50").Value = "55.14"
$ws.Range("E50").Value = "  +7.06%  "
$ws.Range("D51").Value = "73.95"
$ws.Range("E51").Value = "  +3.90%  "
